$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(997, 5).Value = "100M"
$ws.Cells.Item(998, 5).Value = "1pcs"
$ws.Cells.Item(999, 3).Value = "11157407"
Write-Host "E997:" $ws.Cells.Item(997,5).Value()
Write-Host "E998:" $ws.Cells.Item(998,5).Value()
Write-Host "C999:" $ws.Cells.Item(999,3).Value()
